$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the Predicted_ROI (Q) values to two decimals (row 6 rounds to a whole
# number, matching the source data). The Predicted_revenue (R) column holds
# =Q*L + L formulas, so it recalculates automatically once Q changes.
$ws.Range("Q2").Value = 1.19
$ws.Range("Q3").Value = 2.16
$ws.Range("Q4").Value = 1.33
$ws.Range("Q5").Value = 3.14
$ws.Range("Q6").Value = 2
$ws.Range("Q7").Value = 2.03
$ws.Range("Q8").Value = 1.21
$ws.Range("Q9").Value = 1.86

# Display the predicted revenue as a whole number (adds a new "0" numeric
# format cell style used by R2:R9).
$ws.Range("R2:R9").NumberFormat = "0"

# Give the new Predicted_revenue column a readable custom width
# (~13.66 characters, closest this engine's pixel-snapped width model
# can reach to the source file's 13.6640625).
$ws.Range("R1").EntireColumn.ColumnWidth = 12.79

# Restore the sheet's last active-cell selection.
$ws.Range("Q13").Select()
